# feat: add 2022-Q3 data
#
# 1. Duplicate the existing "2022-Q2" sheet (same column layout/styling as the
#    new quarter needs) and place the copy right after "总计", then rename it
#    to "2022-Q3" and overwrite its data with the new quarter's figures.
# 2. Insert the new quarter's summary row at the top of the "总计" data table
#    and push the existing rows down by one (the last existing row -
#    "2020-Q4" - gets duplicated down into the newly created row 8).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- Step 1: create the new "2022-Q3" sheet -------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Fund code / name stay identical to the other quarters for this holding;
# only the size / position figures and the ranking change.
$q3Sheet.Range("D2:G2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "3.05"
$q3Sheet.Range("E2").Value = "97.26"
$q3Sheet.Range("F2").Value = "4.21"
$q3Sheet.Range("G2").Value = "0.1284"
$q3Sheet.Range("H2").Value = 1

# --- Step 2: shift the "总计" rows down and insert the new quarter --------
$totalSheet.Range("A8").Value = 6
$totalSheet.Range("B8").Value = "2020-Q4"
$totalSheet.Range("C8").Value = 2
$totalSheet.Range("D8").Value = 0.87

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q1"
$totalSheet.Range("C7").Value = 2
$totalSheet.Range("D7").Value = 1.22

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q2"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.83

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q3"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.66

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.38

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.4

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.13

# new A8 falls outside the old A1:D7 used range, so it doesn't inherit the
# bordered/centered style the rest of column A uses - copy it over explicitly.
$totalSheet.Range("A7").Copy()
$totalSheet.Range("A8").PasteSpecial(-4122)

# Restore the originally-active sheet ("2020-Q4", now the last tab).
$wb.Worksheets.Item("2020-Q4").Activate()
